$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the value in column B (rows 2-19) by 1
for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value2 = $cell.Value2 + 1
}
